$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 47: normal entry row (copy formatting from row 40, which has the same
#     plain/no-special-height layout: A s=1 date, B s=2 wrapped text, C plain text) ---
$ws.Range("A40:C40").Copy($ws.Range("A47:C47"))
$ws.Range("A47").Value = 43175
$ws.Range("B47").Value = "J'ai commencé à implémenter mon QR Code avec l'aide de mon natel pour faire les tests"
$ws.Range("C47").Value = "3 périodes"

# --- Row 48: new "week" section header (copy formatting from row 41, the
#     preceding section-header row, merged A:C, style s=3) ---
$ws.Range("A41:C41").Copy($ws.Range("A48:C48"))
$ws.Range("A48").Value = "7ème semaine "

# --- Row 50: normal entry row (copy formatting from row 40 again). Filled in
#     before row 49 so the shared-string table gets the same insertion order
#     as the authored workbook (index 50 = this text, index 51 = row 49's). ---
$ws.Range("A40:C40").Copy($ws.Range("A50:C50"))
$ws.Range("A50").Value = 43179
$ws.Range("B50").Value = "J'ai rempli ma documentation. J'ai fini la partie où je décris mes tables dans mon MLD"
$ws.Range("C50").Value = "1 période"

# --- Row 49: entry row with row height 75 (copy formatting from row 42, which
#     already has ht=75) ---
$ws.Range("A42:C42").Copy($ws.Range("A49:C49"))
$ws.Rows.Item(49).RowHeight = 75
$ws.Range("A49").Value = 43179
$ws.Range("B49").Value = "J'ai trouvé un tutoriel, pour scanner les qr code, que j'ai suivi sur youtube à cette adresse : https://www.youtube.com/watch?v=aa0abyOBa28. J'ai fait le tuto et je l'ai testé sur mon natel et ça fonctionne il me ressort des infos du qr code que je scanne. J'ai fini ma fonction qui me permet de recevoir la quantité d'article que j'ai dans ma base de données en scannant le QR Code avec l'id d'un article."
$ws.Range("C49").Value = "4 périodes"

# --- Row 51: trailing empty wrapped-text cell in column B only ---
$ws.Range("B40").Copy($ws.Range("B51"))
$ws.Range("B51").Value = ""

# Update the visible selection to match where the workbook was left (C50).
$ws.Range("C50").Select()
